$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 7 (pairwise_merge) - its "response_collected" value gets folded
# into column F of the row above (row 6) before the row is removed.
$ws.Range("F6").Value = $ws.Range("E7").Value2
$ws.Rows("7:7").Delete()

# After the above deletion, the Wali resub row (originally row 10) is now row 9.
# It simply gains a new "Further_process" value in column F (no row removal here).
$ws.Range("F9").Value = "master_all_responses_Jun-19-2023_to_Jul-14-2023_Wali.csv"

# The SB "video_no_play" row (originally row 17) is now row 16. Its
# "submission_file_name" value is dropped, folded instead into column F of the
# row above it (now row 15, originally row 16), then the row is removed.
$ws.Range("F15").Value = "master_all_responses_SB_Oct-01-2023_to_SB_resub_Oct-01-2023_Sarah.csv"
$ws.Rows("16:16").Delete()

$ws.Range("D14").Select()
